# Commit: "Přidání složky s php + čas na hlavní stránku" — the spreadsheet
# part of this commit is the author ticking off the task in row 23
# ("Na úvodní stránku umístěte funkci PHP pro zápis aktuálního data" /
# put the PHP "current date" function on the main page) as done, switching
# its status cell from "ne" to "ano". Every other change in the diff
# (J4/J6/J7 totals, F23, the new "ano" shared string, the chart cache) is a
# downstream ripple of this one input, so we only need to write the cell;
# the workbook's own formulas recompute the rest.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E23").Value = "ano"

# Match the author's on-screen context at save time: scrolled down so row 8
# is at the top of the view, with E23 as the selected/active cell.
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("E23").Select()
